# Aston Villa Stats — append a second season's row of data ("Append Data")
#
# Summary of the edit being reproduced:
#   1. Insert a new first column ("Season") in front of the existing table,
#      shifting the whole Table1 data from A:M to B:N.
#   2. Fill in the new Season column: "Season" header, "24/25" for the
#      existing row, and a brand-new "23/24" row with a second season's
#      worth of stats appended underneath.
#   3. Resize Table1 (the ListObject) so it covers the new B1:N3 range.
#   4. Remove the old hyperlink that was sitting on the "Top Scorer" cell
#      (L2, formerly K2) so it becomes plain shared-string text, and drop
#      the now-unused "Hyperlink" cell style.
#   5. Apply a plain AutoFilter on the new A1 "Season" header cell (this is
#      what stamps the hidden _FilterDatabase defined name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert new column A; old A:M becomes B:N -----------------------
$ws.Columns.Item(1).Insert()

# --- 2. Write the "Season" column data ----------------------------------
$ws.Range("A1").Value = "Season"
$ws.Range("A2").Value = "24/25"
$ws.Range("A3").Value = "23/24"

# --- Append the 23/24 season row (shifted columns B..N) ----------------
$ws.Range("B3").Value = "Aston Villa "
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 68
$ws.Range("E3").Value = 38
$ws.Range("F3").Value = 20
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 8
$ws.Range("I3").Value = 76
$ws.Range("J3").Value = 61
$ws.Range("K3").Value = 15
$ws.Range("L3").Value = "Ollie Watkins"
$ws.Range("L3").Style = $ws.Range("L2").Style
$ws.Range("M3").Value = 2
$ws.Range("N3").Value = 94

# --- 3. Resize Table1 to cover the new data range -----------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B1:N3"))

# --- 4. Drop the hyperlink on the (now) L2 "Top Scorer" cell -----------
$ws.Range("L2").Hyperlinks.Delete()
$wb.Styles.Item("Hyperlink").Delete()

# --- 5. AutoFilter the new Season header + hidden _FilterDatabase name -
$ws.Range("A1").AutoFilter()
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "='Aston Villa Stats'!`$A`$1:`$A`$1")
$fdb.Visible = $false
